# Collections.pptx — "Anonymous methods and lambda expressions"
#
# 1. Append three new "Title and Content" slides (Delegates / Anonymous
#    Methods / Lambda Expressions) after the existing last slide.
# 2. Bump the cached datetimeFigureOut field (06-Mar-24 -> 07-Mar-24) on
#    the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Set-BulletArial($rng) {
    $b = $rng.ParagraphFormat.Bullet
    $b.Font.Name = "Arial"
    $b.Character = 8226
    $b.Visible = -1
}

function Add-ContentSlide($index, $titleText, $bullets) {
    # $bullets is an array of hashtables: @{ Text = "..."; Level = 1; Size = 20 }
    $layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"
    $slide = $p.Slides.AddSlide($index, $layout)

    # --- Title -----------------------------------------------------
    $title = $slide.Shapes.Title
    $ttr = $title.TextFrame.TextRange
    $ttr.Text = $titleText
    $ttr.Font.Bold = -1
    $ttr.Font.Color.RGB = 255          # BGR -> FF0000 (red)
    $ttr.Font.NameFarEast = "Cambria"

    # --- Body --------------------------------------------------------
    $body = $slide.Shapes.Placeholders.Item(2)
    $tf = $body.TextFrame
    $tf.AutoSize = 2                    # ppAutoSizeTextToFitShape -> normAutofit
    $tr = $tf.TextRange

    # Build every paragraph first (keeps the "lang" stamp on each run),
    # then go back and apply per-paragraph formatting via Characters().
    $first = $bullets[0]
    $tr.Text = $first.Text
    for ($i = 1; $i -lt $bullets.Count; $i++) {
        $tr.InsertAfter("`r" + $bullets[$i].Text) | Out-Null
    }

    $pos = 1
    foreach ($bullet in $bullets) {
        $len = $bullet.Text.Length
        if ($len -eq 0) {
            $rng = $tr.Characters($pos, 1)
        } else {
            $rng = $tr.Characters($pos, $len)
        }
        $rng.IndentLevel = $bullet.Level + 1
        $rng.Font.Size = $bullet.Size
        Set-BulletArial $rng
        $pos += $len + 1
    }

    return $slide
}

# ---------------------------------------------------------------------
# Slide 23 — Delegates (generic delegates)
# ---------------------------------------------------------------------
Add-ContentSlide 23 "Delegates" @(
    @{ Text = "Generic delegates"; Level = 1; Size = 20 },
    @{ Text = "If we have five methods, we can create five delegates. But when we work on large projects, we may have to create fifty or five hundred delegates. If we want to encapsulate them by using delegates, then we must create fifty or five hundred delegates. That will really increase the size of our code and make our application slow."; Level = 2; Size = 16 },
    @{ Text = "This ca be addressed with generic delegates"; Level = 2; Size = 16 },
    @{ Text = ""; Level = 1; Size = 16 }
) | Out-Null

# ---------------------------------------------------------------------
# Slide 24 — Anonymous Methods
# ---------------------------------------------------------------------
Add-ContentSlide 24 "Anonymous Methods" @(
    @{ Text = "The anonymous methods are defined using the delegate keyword"; Level = 1; Size = 20 },
    @{ Text = "An anonymous method must be assigned to a delegate type."; Level = 1; Size = 20 },
    @{ Text = "This method can access outer variables or functions except for the outer function ref and out parameter."; Level = 1; Size = 20 },
    @{ Text = "An anonymous method can be passed as a parameter."; Level = 1; Size = 20 }
) | Out-Null

# ---------------------------------------------------------------------
# Slide 25 — Lambda Expressions
# ---------------------------------------------------------------------
Add-ContentSlide 25 "Lambda Expressions" @(
    @{ Text = "The Lambda Expression in C# is the shorthand for writing the Anonymous Function. "; Level = 1; Size = 20 },
    @{ Text = "So, we can say that the Lambda Expression is nothing but to simplify the anonymous function in C# and,"; Level = 1; Size = 20 },
    @{ Text = "we also discussed that Anonymous Functions are related to delegate and they are created by using the delegate keyword"; Level = 1; Size = 20 }
) | Out-Null

# ---------------------------------------------------------------------
# Bump the cached date field: 06-Mar-24 -> 07-Mar-24 everywhere it is
# cached (slide master + all slide layouts).
# ---------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "06-Mar-24") {
                $tr.Text = "07-Mar-24"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}
